# ToDo list update for No2.4.6.7 - add "Comments" column, mark several
# items as "Done", and record reviewer comments for two of the items.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet1: new "Comments" column (C) ------------------------------------

$ws1.Range("C2").Value = "Comments"
$ws1.Range("C2").HorizontalAlignment = -4108   # xlCenter

# Pre-format the column body (left aligned) down to row 93, matching the
# extent the author pre-formatted in the real workbook.
$ws1.Range("C3:C93").HorizontalAlignment = -4131   # xlLeft

$ws1.Columns.Item(3).ColumnWidth = 50.5

# Mark items 2, 4, 6 and 7 as Done, and leave review comments on items 6 and 7.
$ws1.Range("B4").Value = "Done"
$ws1.Range("B6").Value = "Done"
$ws1.Range("B8").Value = "Done"
$ws1.Range("B9").Value = "Done"

$ws1.Range("C8").Value = "左侧预留出空间，并在显示的时候左侧也预留两个空格"
$ws1.Range("C9").Value = "整行加了个灰色背景"

# --- Sheet2 / Sheet3: bring selection state in line with the rest of the book --

$ws2.Range("A2").Select() | Out-Null

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1   # xlPortrait

# Leave Sheet1 as the active sheet/selection, matching the saved view state.
$ws1.Range("B7").Select() | Out-Null
